$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.231.21'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.011.35'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''246.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '''0.644'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.80%  '
$ws.Range("D7").Value = '''63.15'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +19.37%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '''59.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").Value = '''0.0748'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '''0.947'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = '''14.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.21%  '
$ws.Range("D15").Value = '2.301.61'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '''5.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = '''19.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +17.07%  '
$ws.Range("D18").Value = '2.009.99'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").Value = '36.150.41'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '''72.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").Value = '0.0₃0859'
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").Value = '''5.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.04%  '
$ws.Range("D23").Value = '''234.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '''2.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +18.79%  '
$ws.Range("D26").Value = '''2.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").Value = '''9.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.27%  '
$ws.Range("D28").Value = '''166.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("D29").Value = '''19.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '''0.120'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").Value = '''5.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.07%  '
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '''0.0996'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +16.88%  '
$ws.Range("E34").Value = '  +3.84%  '
$ws.Range("D35").Value = '''4.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.87%  '
$ws.Range("D36").Value = '''2.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.93%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").Value = '''5.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +17.92%  '
$ws.Range("D40").Value = '''1.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.96%  '
$ws.Range("D41").Value = '''0.0964'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.85%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.89'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '''0.0216'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.49%  '
$ws.Range("D44").Value = '''16.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.95%  '
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("D46").Value = '''94.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("D47").Value = '''7.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.69%  '
$ws.Range("D48").Value = '1.373.67'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").Value = '''2.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.08%  '
$ws.Range("D51").Value = '''47.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.82%  '
